$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.080.05"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "'2.639.66"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'596.12"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").Value = "'156.23"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.543"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.142"
$ws.Range("E9").Value = "  +1.82%  "
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("D11").Value = "'5.23"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'0.351"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'27.99"
$ws.Range("E13").Value = "  +0.18%  "
$ws.Range("D14").Value = "'0.0000189"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "'3.121.59"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "'67.939.63"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'2.641.58"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'11.33"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "'362.15"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "'7.39"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "'4.41"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "'4.78"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").Value = "'2.05"
$ws.Range("E23").Value = "  -3.55%  "
$ws.Range("D24").Value = "'75.00"
$ws.Range("E24").Value = "  +2.63%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'9.68"
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").Value = "'552.93"
$ws.Range("E30").Value = "  -5.95%  "
$ws.Range("D31").Value = "'7.98"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "'1.84"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("D36").Value = "'1.54"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").Value = "'161.17"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "'19.36"
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'0.371"
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").Value = "'5.30"
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("D42").Value = "'0.0₆0338"
$ws.Range("E42").Value = "  +5.88%  "
$ws.Range("D44").Value = "'2.60"
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "'158.85"
$ws.Range("E47").Value = "  +1.37%  "
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").Value = "'21.97"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0784"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D51").Value = "'1.69"
$ws.Range("E51").Value = "  -1.81%  "
